# Update the "ActlHours" and "ActDeliveryDate" for the task in row 19
# (Inloggningssida / Funktion) to reflect the continued work logged on
# the user authentication endpoint.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Actual hours spent increased from 32 to 43
$ws.Range("P19").Value = 43

# The pending delivery date moved from 2020-02-28 to 2020-03-03; this is
# stored as plain text (not a real date) in the sheet.
$ws.Range("Q19").Value = "2020-03-03 Pending"

$wb.Save()
